# Apply crypto price/volume updates as described by the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '54.596.60'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.54%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.293.43'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -4.28%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.997'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '496.14'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '127.42'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -4.85%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.531'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.291.66'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -4.22%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.50%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.325'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -4.92%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.685.75'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -4.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.65'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.62%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '54.496.97'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.64%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.298.47'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -5.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.04'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.38%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '304.73'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.38'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.59%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.26'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.78%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.90%  '
$ws.Range("B27").Value = 'Polygon'
$ws.Range("C27").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.374'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.04%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.151'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.56%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.364.09'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -4.85%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.73%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '170.28'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.72%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0₃0687'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -4.95%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.23%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.24%  '
$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.25%  '
$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.08'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.60%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '17.60'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.18%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.863'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.29%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '35.55'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.77%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.375'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.66%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '129.86'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.14%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.88%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.82'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -4.56%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.89%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.547'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.87%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '241.90'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.44%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.40%  '
